# Scheduled runner update: refresh Market Board price snapshots (currentAveragePrice*,
# LevePrice*, LeveProfit*) across the ARM/BSM/CRP/CUL/GSM/LTW/WVR profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5493.0347
$ws.Range("I32").Value = 5200.393
$ws.Range("J32").Value = 17784
$ws.Range("K32").Value = 5200.393
$ws.Range("L32").Value = 17784
$ws.Range("M32").Value = -4913.393
$ws.Range("N32").Value = -18358
$ws.Range("H56").Value = 10600
$ws.Range("J56").Value = 10600
$ws.Range("L56").Value = 10600
$ws.Range("N56").Value = -12084
$ws.Range("H61").Value = 1652.8055
$ws.Range("I61").Value = 1185.25
$ws.Range("J61").Value = 1886.5834
$ws.Range("K61").Value = 1185.25
$ws.Range("L61").Value = 1886.5834
$ws.Range("M61").Value = -973.25
$ws.Range("N61").Value = -2310.5834
$ws.Range("H109").Value = 16600
$ws.Range("J109").Value = 16600
$ws.Range("L109").Value = 16600
$ws.Range("N109").Value = -19374
$ws.Range("H136").Value = 1652.8055
$ws.Range("I136").Value = 1185.25
$ws.Range("J136").Value = 1886.5834
$ws.Range("K136").Value = 3555.75
$ws.Range("L136").Value = 5659.7502
$ws.Range("M136").Value = -1005.75
$ws.Range("N136").Value = -10759.7502

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 15000
$ws.Range("I8").Value = 15000
$ws.Range("K8").Value = 15000
$ws.Range("M8").Value = -14860
$ws.Range("H10").Value = 37668.668
$ws.Range("I10").Value = 10000
$ws.Range("J10").Value = 51503
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 51503
$ws.Range("M10").Value = -9860
$ws.Range("N10").Value = -51783
$ws.Range("H108").Value = 40000
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680
$ws.Range("H109").Value = 28170.5
$ws.Range("J109").Value = 28170.5
$ws.Range("L109").Value = 28170.5
$ws.Range("N109").Value = -30944.5
$ws.Range("H112").Value = 15071.167
$ws.Range("J112").Value = 15071.167
$ws.Range("L112").Value = 15071.167
$ws.Range("N112").Value = -18025.167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 17551.625
$ws.Range("J43").Value = 17551.625
$ws.Range("L43").Value = 17551.625
$ws.Range("N43").Value = -17919.625
$ws.Range("H60").Value = 9198.6
$ws.Range("I60").Value = 3998.25
$ws.Range("K60").Value = 3998.25
$ws.Range("M60").Value = -3487.25
$ws.Range("H101").Value = 17551.625
$ws.Range("J101").Value = 17551.625
$ws.Range("L101").Value = 17551.625
$ws.Range("N101").Value = -24041.625
$ws.Range("H132").Value = 2624.3
$ws.Range("I132").Value = 2172.7
$ws.Range("J132").Value = 2850.1
$ws.Range("K132").Value = 6518.099999999999
$ws.Range("L132").Value = 8550.299999999999
$ws.Range("M132").Value = -3988.099999999999
$ws.Range("N132").Value = -13610.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 1860
$ws.Range("I58").Value = 780
$ws.Range("J58").Value = 2400
$ws.Range("K58").Value = 2340
$ws.Range("L58").Value = 7200
$ws.Range("M58").Value = -2212
$ws.Range("N58").Value = -7456
$ws.Range("H70").Value = 3746.8462
$ws.Range("I70").Value = 1542.1428
$ws.Range("J70").Value = 6319
$ws.Range("K70").Value = 4626.428400000001
$ws.Range("L70").Value = 18957
$ws.Range("M70").Value = -4311.428400000001
$ws.Range("N70").Value = -19587
$ws.Range("H73").Value = 3746.8462
$ws.Range("I73").Value = 1542.1428
$ws.Range("J73").Value = 6319
$ws.Range("K73").Value = 4626.428400000001
$ws.Range("L73").Value = 18957
$ws.Range("M73").Value = -3534.428400000001
$ws.Range("N73").Value = -21141
$ws.Range("H122").Value = 1712.9849
$ws.Range("I122").Value = 554.0769
$ws.Range("J122").Value = 1997.2452
$ws.Range("K122").Value = 4986.6921
$ws.Range("L122").Value = 17975.2068
$ws.Range("M122").Value = -2536.6921
$ws.Range("N122").Value = -22875.2068

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6228.1665
$ws.Range("I70").Value = 7000.727
$ws.Range("J70").Value = 5014.143
$ws.Range("K70").Value = 7000.727
$ws.Range("L70").Value = 5014.143
$ws.Range("M70").Value = -6730.727
$ws.Range("N70").Value = -5554.143
$ws.Range("H73").Value = 6228.1665
$ws.Range("I73").Value = 7000.727
$ws.Range("J73").Value = 5014.143
$ws.Range("K73").Value = 7000.727
$ws.Range("L73").Value = 5014.143
$ws.Range("M73").Value = -6064.727
$ws.Range("N73").Value = -6886.143
$ws.Range("H97").Value = 1150.75
$ws.Range("I97").Value = 728.4286
$ws.Range("J97").Value = 1742
$ws.Range("K97").Value = 728.4286
$ws.Range("L97").Value = 728.4286
$ws.Range("M97").Value = -232.4286
$ws.Range("N97").Value = -2734

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 11551
$ws.Range("J104").Value = 11551
$ws.Range("L104").Value = 11551
$ws.Range("N104").Value = -18539
$ws.Range("H111").Value = 30795.666
$ws.Range("J111").Value = 30795.666
$ws.Range("L111").Value = 30795.666
$ws.Range("N111").Value = -38975.666
$ws.Range("H132").Value = 33761.383
$ws.Range("I132").Value = 54466.2
$ws.Range("J132").Value = 4183.0713
$ws.Range("K132").Value = 163398.6
$ws.Range("L132").Value = 12549.2139
$ws.Range("M132").Value = -160868.6
$ws.Range("N132").Value = -17609.2139
$ws.Range("H136").Value = 1687.6154
$ws.Range("I136").Value = 1472.68
$ws.Range("J136").Value = 2071.4285
$ws.Range("K136").Value = 4418.04
$ws.Range("L136").Value = 6214.2855
$ws.Range("M136").Value = -1868.04
$ws.Range("N136").Value = -11314.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3278.1428
$ws.Range("I62").Value = 2822.6
$ws.Range("J62").Value = 4417
$ws.Range("K62").Value = 2822.6
$ws.Range("L62").Value = 4417
$ws.Range("M62").Value = -2198.6
$ws.Range("N62").Value = -5665
$ws.Range("H65").Value = 3278.1428
$ws.Range("I65").Value = 2822.6
$ws.Range("J65").Value = 4417
$ws.Range("K65").Value = 14113
$ws.Range("L65").Value = 22085
$ws.Range("M65").Value = -10993
$ws.Range("N65").Value = -28325
$ws.Range("H69").Value = 8269.916999999999
$ws.Range("J69").Value = 8269.916999999999
$ws.Range("L69").Value = 8269.916999999999
$ws.Range("N69").Value = -9767.916999999999
$ws.Range("H72").Value = 8269.916999999999
$ws.Range("J72").Value = 8269.916999999999
$ws.Range("L72").Value = 24809.751
$ws.Range("N72").Value = -32297.751
$ws.Range("H82").Value = 9994.5
$ws.Range("J82").Value = 9994.5
$ws.Range("L82").Value = 9994.5
$ws.Range("N82").Value = -10760.5
$ws.Range("H85").Value = 9994.5
$ws.Range("J85").Value = 9994.5
$ws.Range("L85").Value = 9994.5
$ws.Range("N85").Value = -12646.5
$ws.Range("H126").Value = 9215.846
$ws.Range("I126").Value = 9089.556
$ws.Range("J126").Value = 9500
$ws.Range("K126").Value = 27268.668
$ws.Range("L126").Value = 28500
$ws.Range("M126").Value = -24798.668
$ws.Range("N126").Value = -33440
$ws.Range("H132").Value = 1726.491
$ws.Range("I132").Value = 1457.6061
$ws.Range("J132").Value = 2129.818
$ws.Range("K132").Value = 4372.8183
$ws.Range("L132").Value = 6389.454000000001
$ws.Range("M132").Value = -1842.8183
$ws.Range("N132").Value = -11449.454
$ws.Range("H136").Value = 1350.8
$ws.Range("I136").Value = 1180.4131
$ws.Range("K136").Value = 3541.2393
$ws.Range("M136").Value = -991.2393000000002
